# Update "想去人数" (want-to-go count) figures in column F across the
# "展览" (Exhibitions), "演出" (Performances) and "全部类型" (All types)
# worksheets, matching the refreshed snapshot of the source data feed.

$wb = $excel.ActiveWorkbook

$sheetExhibitions = $wb.Worksheets.Item("展览")
$sheetPerformances = $wb.Worksheets.Item("演出")
$sheetAll = $wb.Worksheets.Item("全部类型")

# 展览 (Exhibitions) sheet updates
$sheetExhibitions.Range("F2").Value = 1324
$sheetExhibitions.Range("F4").Value = 79
$sheetExhibitions.Range("F9").Value = 1007
$sheetExhibitions.Range("F11").Value = 175
$sheetExhibitions.Range("F15").Value = 345
$sheetExhibitions.Range("F16").Value = 753
$sheetExhibitions.Range("F17").Value = 135
$sheetExhibitions.Range("F19").Value = 253
$sheetExhibitions.Range("F21").Value = 972
$sheetExhibitions.Range("F22").Value = 433
$sheetExhibitions.Range("F23").Value = 242
$sheetExhibitions.Range("F28").Value = 454

# 演出 (Performances) sheet updates
$sheetPerformances.Range("F6").Value = 39

# 全部类型 (All types) sheet updates
$sheetAll.Range("F3").Value = 1324
$sheetAll.Range("F6").Value = 79
$sheetAll.Range("F11").Value = 1007
$sheetAll.Range("F13").Value = 175
$sheetAll.Range("F20").Value = 39
$sheetAll.Range("F22").Value = 345
$sheetAll.Range("F23").Value = 753
$sheetAll.Range("F24").Value = 135
$sheetAll.Range("F26").Value = 253
$sheetAll.Range("F28").Value = 972
$sheetAll.Range("F29").Value = 433
$sheetAll.Range("F32").Value = 242
$sheetAll.Range("F40").Value = 454
